$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank row 36 (B: date, C: hours, D: description)
$ws.Range("B36").Value = "5/4/2024"
$ws.Range("C36").Value = 4.5
$ws.Range("D36").Value = "Configured correctly Oauth2 flow"

# Update the selection to match the author's final selection (D39:D41)
$ws.Range("D39:D41").Select()
